$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Lecture 2.x topic names to reflect new numbering (2.1 / 2.2)
$ws.Range("C3").Value = "2.1 - Circuits and Layout Part 1"
$ws.Range("C4").Value = "2.2 - Circuits and Layout Part 2"

# Shift the "Due" column (G) assignments: HW 1 / Quiz 1 / Lab 1 each move
# down one row, freeing G10 and appending a new entry at G13.
$ws.Range("G10").ClearContents()
$ws.Range("G11").Value = "HW 1"
$ws.Range("G12").Value = "Quiz 1"
$ws.Range("G13").Value = "Lab 1"

# Add the Lecture 5.0 video recording link in H9
$ws.Range("H9").Value = "https://iu.zoom.us/rec/share/Xubuchvq8ycmR_WQz2GwuMwlpqEnvuU-Nph_SxxNFBwEU1oLzMLGcA_NrpX54rg.UX0vVf-SQbjjsG-A"
$ws.Hyperlinks.Add($ws.Range("H9"), "https://iu.zoom.us/rec/share/Xubuchvq8ycmR_WQz2GwuMwlpqEnvuU-Nph_SxxNFBwEU1oLzMLGcA_NrpX54rg.UX0vVf-SQbjjsG-A")
$ws.Range("H9").Style = $ws.Range("H8").Style

# Update the saved selection to match the author's final cursor position
$ws.Range("F18").Select()
